$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("P2").Value = 0.01
$ws.Range("Q2").Value = -1.086733406199907
$ws.Range("R2").Value = 6.232284007127456
$ws.Range("S2").Value = -5.77383943227396

# Update row 3 values
$ws.Range("Q3").Value = -1.193628943637368
$ws.Range("R3").Value = 9.526209573536834
$ws.Range("S3").Value = -5.786055553696495

# Update the label strings referenced by U2/U3 (shared strings content changes
# from "blink"/"templerun" to "blink+templerun"/"blink+sudoku")
$ws.Range("U2").Value = "blink+templerun"
$ws.Range("U3").Value = "blink+sudoku"

# Delete rows 4 through 6 entirely, shrinking the used range to A1:U3
$ws.Range("A4:U6").EntireRow.Delete()
